$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "23_00" sheet as a copy of the last existing sheet
#    ("22_45"), placed after it, and make it the selected tab.
# ---------------------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("22_45")
$sourceSheet.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "23_00"

# ---------------------------------------------------------------------------
# 2. Header block (rows 1-6)
# ---------------------------------------------------------------------------

# Row 1: "Apylinkių skaičius - 1895  (duomenys iš 1172 apylinkių)" + count cell
$a1 = $ws.Range("A1")
$a1full = "Apylinkių skaičius - 1895  (duomenys iš 1172 apylinkių)"
$a1.Value = $a1full
$r = $a1.Characters(22, 4)           # "1895"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a1.Characters(26, 1)           # " "
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a1.Characters(27, 29)          # " (duomenys iš 1172 apylinkių)"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r.Font.Color = 255

$ws.Range("B1").Value = 1172

# Row 4: " rinkėjų sąraše įrašyta rinkėjų - 752758, rinkimuose dalyvavo - 420889 (55,91%),"
$a4 = $ws.Range("A4")
$a4full = " rinkėjų sąraše įrašyta rinkėjų - 752758, rinkimuose dalyvavo - 420889 (55,91%),"
$a4.Value = $a4full
$r = $a4.Characters(35, 6)           # "752758"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a4.Characters(41, 24)          # ", rinkimuose dalyvavo - "
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a4.Characters(65, 6)           # "420889"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a4.Characters(71, 2)           # " ("
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a4.Characters(73, 5)           # "55,91"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a4.Characters(78, 3)           # "%),"
$r.Font.Name = "Arial"
$r.Font.Size = 11

# Row 5: " negaliojančių biuletenių - 3444 (0,82%), galiojančių biuletenių - 417445 (99,18%)."
$a5 = $ws.Range("A5")
$a5full = " negaliojančių biuletenių - 3444 (0,82%), galiojančių biuletenių - 417445 (99,18%)."
$a5.Value = $a5full
$r = $a5.Characters(29, 4)           # "3444"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(33, 2)           # " ("
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(35, 4)           # "0,82"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(39, 29)          # "%), galiojančių biuletenių - "
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(68, 6)           # "417445"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(74, 2)           # " ("
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(76, 5)           # "99,18"
$r.Font.Bold = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r = $a5.Characters(81, 3)           # "%)."
$r.Font.Name = "Arial"
$r.Font.Size = 11

# ---------------------------------------------------------------------------
# 3. Results table (rows 11-19) with the new vote counts. Candidates 2 and 3
#    (Šimonytė / Žemaitaitis) swap places relative to the previous sheet.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Gitanas NAUSĖDA"
$ws.Range("B11").Value = 155624
$ws.Range("C11").Value = 49605
$ws.Range("D11").Value = 205229
$ws.Range("E11").Value = 49.16
$ws.Range("F11").Value = 48.76

$ws.Range("A12").Value = "Ingrida ŠIMONYTĖ"
$ws.Range("B12").Value = 38245
$ws.Range("C12").Value = 15414
$ws.Range("D12").Value = 53659
$ws.Range("E12").Value = 12.85
$ws.Range("F12").Value = 12.75

$ws.Range("A13").Value = "Remigijus ŽEMAITAITIS"
$ws.Range("B13").Value = 44999
$ws.Range("C13").Value = 7349
$ws.Range("D13").Value = 52348
$ws.Range("E13").Value = 12.54
$ws.Range("F13").Value = 12.44

$ws.Range("A14").Value = "Ignas VĖGĖLĖ"
$ws.Range("B14").Value = 42023
$ws.Range("C14").Value = 9051
$ws.Range("D14").Value = 51074
$ws.Range("E14").Value = 12.23
$ws.Range("F14").Value = 12.13

$ws.Range("A15").Value = "Eduardas VAITKUS"
$ws.Range("B15").Value = 27006
$ws.Range("C15").Value = 4968
$ws.Range("D15").Value = 31974
$ws.Range("E15").Value = 7.66
$ws.Range("F15").Value = 7.6

$ws.Range("A16").Value = "Dainius ŽALIMAS"
$ws.Range("B16").Value = 7882
$ws.Range("C16").Value = 3047
$ws.Range("D16").Value = 10929
$ws.Range("E16").Value = 2.62
$ws.Range("F16").Value = 2.6

$ws.Range("A17").Value = "Andrius MAZURONIS"
$ws.Range("B17").Value = 4994
$ws.Range("C17").Value = 1845
$ws.Range("D17").Value = 6839
$ws.Range("E17").Value = 1.64
$ws.Range("F17").Value = 1.62

$ws.Range("A18").Value = "Giedrimas JEGLINSKAS"
$ws.Range("B18").Value = 3910
$ws.Range("C18").Value = 1483
$ws.Range("D18").Value = 5393
$ws.Range("E18").Value = 1.29
$ws.Range("F18").Value = 1.28

$ws.Range("A19").Value = "Iš viso"
$ws.Range("B19").Value = 324683
$ws.Range("C19").Value = 92762
$ws.Range("D19").Value = 417445
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 99.18

# ---------------------------------------------------------------------------
# 4. Hyperlinks: each candidate keeps the same target URL, so A12/A13 need
#    their links swapped along with the row contents above.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @(
    @{ Cell = "A11"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435326.html" },
    @{ Cell = "A12"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435331.html" },
    @{ Cell = "A13"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435327.html" },
    @{ Cell = "A14"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435334.html" },
    @{ Cell = "A15"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435325.html" },
    @{ Cell = "A16"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435333.html" },
    @{ Cell = "A17"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435329.html" },
    @{ Cell = "A18"; Url = "https://rezultatai.vrk.lt/?srcUrl=/rinkimai/1504/1/2070/rezultatai/lt/rezultataiPreKandBalsLietuvoje_rkndId-2435328.html" }
)

foreach ($link in $links) {
    $cell = $ws.Range($link.Cell)
    $ws.Hyperlinks.Add($cell, $link.Url)
    $cell.Style = "Lien hypertexte"
}
